$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 10
$aw.ScrollColumn = 3
Write-Host "done"
